$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.118.00"
$ws.Range("E2").Value = "  -0.92%  "
$ws.Range("D3").Value = "2.943.81"
$ws.Range("E3").Value = "  -1.58%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "375.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.20"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.64%  "
$ws.Range("E7").Value = "  -1.52%  "
$ws.Range("E9").Value = "  -2.06%  "
$ws.Range("E10").Value = "  -3.03%  "
$ws.Range("E11").Value = "  -0.94%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0850"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.19%  "
$ws.Range("D13").Value = "3.403.80"
$ws.Range("E13").Value = "  -1.18%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.05"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.55"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.46%  "
$ws.Range("D16").Value = "2.954.96"
$ws.Range("E16").Value = "  -0.91%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.994"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.61%  "
$ws.Range("E18").Value = "  +43.44%  "
$ws.Range("D19").Value = "50.970.89"
$ws.Range("E19").Value = "  -1.07%  "
$ws.Range("E20").Value = "  -7.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.43"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.68%  "
$ws.Range("E22").Value = "  -0.74%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "265.87"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.70"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.02%  "
$ws.Range("E25").Value = "  +8.47%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.15"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.66"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.93%  "
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "25.61"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.78%  "
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.164"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.67%  "
$ws.Range("E31").Value = "  -6.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "50.67"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.78%  "
$ws.Range("E34").Value = "  -1.41%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "33.33"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.52%  "
$ws.Range("E36").Value = "  -2.71%  "
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.16"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.05%  "
$ws.Range("E39").Value = "  -1.40%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.24"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.68%  "
$ws.Range("E41").Value = "  -3.71%  "
$ws.Range("E42").Value = "  -4.27%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "120.49"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.35"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.29%  "
$ws.Range("E45").Value = "  -1.16%  "
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.34"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.67%  "
$ws.Range("B47").Value = "TheGraph"
$ws.Range("C47").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.272"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.12%  "
$ws.Range("E48").Value = "  -2.76%  "
$ws.Range("D49").Value = "1.991.57"
$ws.Range("E49").Value = "  -2.45%  "
$ws.Range("E50").Value = "  -2.55%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.25%  "
